$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$fileName = "c8529d5f-d313-4d2b-a8bd-01f7b8d9050a.md"
$pathAndName = "e2e\c8529d5f-d313-4d2b-a8bd-01f7b8d9050a.md"
$extension = ".md"
$status = "Ready for handoff"
$sourcePath = "e2e"
$priority = "ht"
$contentDuplicate = "False"
$toBeLocalized = "True"
$handbackDateTime = "0001-01-01 00:00:00"
$handoffDateTimeZhCn = "2016-08-15 20:41:15"
$handoffDateTimeDeDe = "2016-08-15 20:41:20"
$overviewGenDate = "2016-08-15 20:41:20"
$handoffFileZhCn = "c8529d5f-d313-4d2b-a8bd-01f7b8d9050a.f25970f798478e8e5970158fc45067fbbb2b9345.zh-cn.xlf"
$handoffFileDeDe = "c8529d5f-d313-4d2b-a8bd-01f7b8d9050a.f25970f798478e8e5970158fc45067fbbb2b9345.de-de.xlf"
$mdHyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f25970f798478e8e5970158fc45067fbbb2b9345/e2e/c8529d5f-d313-4d2b-a8bd-01f7b8d9050a.md"

# ---------------------------------------------------------------------------
# Overview sheet (row 9): File Name | Path And Name | Extension | Publish URL
#                          | zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$ws1.Range("A9").Value = $fileName
$ws1.Range("C9").Value = $extension
$ws1.Range("D9").Value = ""
$ws1.Range("E9").Value = $status
$ws1.Range("F9").Value = $status
$ws1.Range("G9").Value = $overviewGenDate
$ws1.Range("G9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Hyperlinks.Add($ws1.Range("B9"), $mdHyperlinkAddress, "", "", $pathAndName)

# ---------------------------------------------------------------------------
# zh-cn sheet (row 9)
# ---------------------------------------------------------------------------
$ws2.Range("B9").Value = $extension
$ws2.Range("C9").Value = $status
$ws2.Range("D9").Value = $sourcePath
$ws2.Range("E9").Value = $priority
$ws2.Range("F9").Value = $contentDuplicate
$ws2.Range("G9").Value = $handoffFileZhCn
$ws2.Range("H9").Value = $handoffDateTimeZhCn
$ws2.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("I9").Value = ""
$ws2.Range("J9").Value = ""
$ws2.Range("K9").Value = $handbackDateTime
$ws2.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("L9").Value = ""
$ws2.Range("M9").Value = $toBeLocalized
$ws2.Range("N9").Value = ""
$ws2.Range("O9").Value = $contentDuplicate
$ws2.Range("P9").Value = ""

$ws2.Hyperlinks.Add($ws2.Range("A9"), $mdHyperlinkAddress, "", "", $fileName)

# ---------------------------------------------------------------------------
# de-de sheet (row 9)
# ---------------------------------------------------------------------------
$ws3.Range("B9").Value = $extension
$ws3.Range("C9").Value = $status
$ws3.Range("D9").Value = $sourcePath
$ws3.Range("E9").Value = $priority
$ws3.Range("F9").Value = $contentDuplicate
$ws3.Range("G9").Value = $handoffFileDeDe
$ws3.Range("H9").Value = $handoffDateTimeDeDe
$ws3.Range("H9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("I9").Value = ""
$ws3.Range("J9").Value = ""
$ws3.Range("K9").Value = $handbackDateTime
$ws3.Range("K9").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("L9").Value = ""
$ws3.Range("M9").Value = $toBeLocalized
$ws3.Range("N9").Value = ""
$ws3.Range("O9").Value = $contentDuplicate
$ws3.Range("P9").Value = ""

$ws3.Hyperlinks.Add($ws3.Range("A9"), $mdHyperlinkAddress, "", "", $fileName)

# ---------------------------------------------------------------------------
# Resize the three tables + autofilter ranges to include the new row
# ---------------------------------------------------------------------------
$ws1.ListObjects.Item(1).Resize($ws1.Range("A1:G9"))
$ws2.ListObjects.Item(1).Resize($ws2.Range("A1:P9"))
$ws3.ListObjects.Item(1).Resize($ws3.Range("A1:P9"))

Write-Host "Applied handback report row for c8529d5f-d313-4d2b-a8bd-01f7b8d9050a.md"
